$d = $word.ActiveDocument

# 1. Heading + bold "CTA" text: both occurrences share identical text, Find/Replace handles all matches.
$d.Content.Find.Execute(
    "Play Elements: The Awakening Online Slot Game for Free", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Elements: The Awakening Free", 2) | Out-Null

# 2. "What we like" bullets - reshuffled/reworded.
$d.Content.Find.Execute(
    "Unique fantasy interpretation of natural elements", $true, $false, $false, $false, $false,
    $true, 1, $false, "Versatile themes and unique games", 2) | Out-Null

$d.Content.Find.Execute(
    "Cascading mechanism for concatenated wins", $true, $false, $false, $false, $false,
    $true, 1, $false, "Fantasy interpretation of natural elements", 2) | Out-Null

$d.Content.Find.Execute(
    "Outstanding visual aspects and sound", $true, $false, $false, $false, $false,
    $true, 1, $false, "Cascading mechanism for exciting gameplay", 2) | Out-Null

$d.Content.Find.Execute(
    "High-quality game developer", $true, $false, $false, $false, $false,
    $true, 1, $false, "Chance for concatenated wins", 2) | Out-Null

# 3. "What we don't like" bullets - first bullet reworded, second bullet removed entirely.
$d.Content.Find.Execute(
    "No progressive jackpot feature", $true, $false, $false, $false, $false,
    $true, 1, $false, "Finding similar slots with the theme may be challenging", 2) | Out-Null

$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Limited slot selection for combining fantasy and natural themes*") {
        $p.Range.Delete()
    }
}

# 4. Closing italic summary paragraph. Set the run text directly (instead of
# Find/Replace) so the straight apostrophe in "NetEnt's" survives instead of
# being smart-quoted, and the paragraph's italic run formatting is preserved.
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Explore the features and gameplay of NetEnt*") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "Unbiased review of NetEnt's Elements: The Awakening online slot game. Play for free and win big."
    }
}
